$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on the price cells whose values would otherwise
# be auto-coerced to numeric by Excel (single-dot decimal-looking strings),
# so they stay literal text exactly like the source data feed produces.
$textCells = @("D5","D6","D7","D9","D10","D11","D13","D14","D17","D19","D20","D21","D23","D24","D25","D26","D29","D30","D31","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '52.235.61'
$ws.Range("E2").Value = '  +0.96%  '

# Row 3
$ws.Range("D3").Value = '2.912.67'
$ws.Range("E3").Value = '  +3.93%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '351.58'
$ws.Range("E5").Value = '  -1.23%  '

# Row 6
$ws.Range("D6").Value = '112.52'
$ws.Range("E6").Value = '  +3.13%  '

# Row 7
$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +0.73%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +0.62%  '

# Row 10
$ws.Range("D10").Value = '40.01'
$ws.Range("E10").Value = '  +0.32%  '

# Row 11
$ws.Range("D11").Value = '0.0862'
$ws.Range("E11").Value = '  +2.74%  '

# Row 12
$ws.Range("E12").Value = '  +0.33%  '

# Row 13
$ws.Range("D13").Value = '20.06'
$ws.Range("E13").Value = '  +0.61%  '

# Row 14
$ws.Range("D14").Value = '7.83'
$ws.Range("E14").Value = '  +0.89%  '

# Row 15
$ws.Range("D15").Value = '3.370.77'
$ws.Range("E15").Value = '  +3.91%  '

# Row 16
$ws.Range("D16").Value = '2.924.47'
$ws.Range("E16").Value = '  +5.01%  '

# Row 17
$ws.Range("D17").Value = '0.999'
$ws.Range("E17").Value = '  +6.08%  '

# Row 18
$ws.Range("D18").Value = '52.279.34'
$ws.Range("E18").Value = '  +1.07%  '

# Row 19
$ws.Range("D19").Value = '7.66'
$ws.Range("E19").Value = '  -0.25%  '

# Row 20
$ws.Range("D20").Value = '3.31'
$ws.Range("E20").Value = '  +5.35%  '

# Row 21
$ws.Range("D21").Value = '14.29'
$ws.Range("E21").Value = '  +5.10%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0982'
$ws.Range("E22").Value = '  +0.27%  '

# Row 23
$ws.Range("D23").Value = '70.94'
$ws.Range("E23").Value = '  +0.75%  '

# Row 24
$ws.Range("D24").Value = '270.14'
$ws.Range("E24").Value = '  +0.75%  '

# Row 25
$ws.Range("D25").Value = '2.80'
$ws.Range("E25").Value = '  +1.54%  '

# Row 26
$ws.Range("D26").Value = '26.79'
$ws.Range("E26").Value = '  +2.68%  '

# Row 27
$ws.Range("E27").Value = '  +0.05%  '

# Row 28
$ws.Range("E28").Value = '  +0.74%  '

# Row 29
$ws.Range("D29").Value = '10.62'
$ws.Range("E29").Value = '  +2.49%  '

# Row 30
$ws.Range("D30").Value = '37.70'
$ws.Range("E30").Value = '  +0.16%  '

# Row 31
$ws.Range("D31").Value = '6.53'
$ws.Range("E31").Value = '  +4.78%  '

# Row 32
$ws.Range("E32").Value = '  +1.29%  '

# Row 33
$ws.Range("E33").Value = '  +8.27%  '

# Row 34
$ws.Range("D34").Value = '0.0956'
$ws.Range("E34").Value = '  +11.25%  '

# Row 35
$ws.Range("D35").Value = '53.32'
$ws.Range("E35").Value = '  +2.72%  '

# Row 36
$ws.Range("E36").Value = '  +1.92%  '

# Row 37
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.02%  '

# Row 38
$ws.Range("D38").Value = '3.31'
$ws.Range("E38").Value = '  +5.04%  '

# Row 39
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '2.08'
$ws.Range("E39").Value = '  +3.91%  '

# Row 40
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '18.77'
$ws.Range("E40").Value = '  -0.36%  '

# Row 41
$ws.Range("D41").Value = '2.80'
$ws.Range("E41").Value = '  +12.72%  '

# Row 42
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").Value = '  +1.63%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '23.55'
$ws.Range("E43").Value = '  +6.85%  '

# Row 44
$ws.Range("D44").Value = '2.65'
$ws.Range("E44").Value = '  +8.10%  '

# Row 45
$ws.Range("D45").Value = '121.66'
$ws.Range("E45").Value = '  +2.12%  '

# Row 46
$ws.Range("E46").Value = '  -0.50%  '

# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.202.44'
$ws.Range("E47").Value = '  +4.49%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '3.53'
$ws.Range("E48").Value = '  +4.17%  '

# Row 49
$ws.Range("D49").Value = '0.263'
$ws.Range("E49").Value = '  +23.84%  '

# Row 50
$ws.Range("D50").Value = '0.0337'
$ws.Range("E50").Value = '  +12.19%  '

# Row 51
$ws.Range("D51").Value = '0.963'
$ws.Range("E51").Value = '  +6.21%  '
